# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and fixes the Algorand / TheSandbox row ordering (rows 37-38).
#
# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (and so lose formatting, e.g. "1.000" -> 1 or "0.05990" -> 0.0599) are
# written with a leading apostrophe so Excel stores them as literal text,
# exactly like the source inline-string values they replace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.590.62"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "1.753.79"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'324.39"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4485"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "'0.07512"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").Value = "'42.19"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("D11").Value = "'1.107"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'20.73"
$ws.Range("E13").Value = "  -5.75%  "
$ws.Range("D14").Value = "'6.051"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").Value = "'7.183"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "1.748.48"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "'92.90"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "'0.06396"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'16.92"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("D22").Value = "'5.859"
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("D23").Value = "27.618.01"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "'11.23"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "'2.106"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "'161.59"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "'20.40"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "1.953.24"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").Value = "'2.124"
$ws.Range("E29").Value = "  -6.38%  "
$ws.Range("D30").Value = "'125.40"
$ws.Range("E30").Value = "  -4.05%  "
$ws.Range("D31").Value = "'1.084"
$ws.Range("E31").Value = "  -10.09%  "
$ws.Range("D32").Value = "'0.09031"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'3.645"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").Value = "'5.564"
$ws.Range("E34").Value = "  -7.68%  "
$ws.Range("D35").Value = "'12.03"
$ws.Range("E35").Value = "  -7.92%  "
$ws.Range("D36").Value = "'0.02316"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "'0.6412"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2096"
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").Value = "'4.993"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").Value = "'0.05990"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'1.197"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.9997"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("D44").Value = "'7.803"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").Value = "'13.31"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "'0.5921"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "'3.713"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'1.961"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("D49").Value = "'121.68"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  -1.73%  "
